$d = $word.ActiveDocument

# Replace "Prezidentė" with "Prezidentas" everywhere in the document (wildcard off)
$d.Content.Find.Execute("Prezidentė", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Prezidentas", 2)

# Replace "Kristė Skaudaitė" with "Danas Černeckas" everywhere in the document
$d.Content.Find.Execute("Kristė Skaudaitė", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Danas Černeckas", 2)
